# Apply row data permutation for rows 2-12 in columns A,B,D,E,F,G,H,Q,R.
# The underlying records (Id, Taxonsorteringsordning, Rödlistade, TaxonId,
# Artnamn, Vetenskapligt namn, Auktor, Ost, Nord) were re-shuffled across
# the existing rows (same set of records, new row order), while columns
# C,I..AY (and everything else) stay attached to their original row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together with each record.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Mapping: new row -> old row that supplies the record's data.
$mapping = @{
    2  = 6
    3  = 2
    4  = 7
    5  = 3
    6  = 8
    7  = 9
    8  = 10
    9  = 11
    10 = 12
    11 = 4
    12 = 5
}

# Snapshot all old values first (since we will overwrite rows in place).
# Use Value2 for reading (Value getter is unreliable for reads in this
# runtime), Value for writing (per the documented usage pattern).
$old = @{}
foreach ($r in 2..12) {
    $old[$r] = @{}
    foreach ($col in $cols) {
        $old[$r][$col] = $ws.Range("$col$r").Value2
    }
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $old[$oldRow][$col]
    }
}
